$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 266, shifting existing rows 266:334 down to 267:335
$ws.Rows("266:266").Insert()

# Populate the newly inserted row 266 with the new weekly data point
$ws.Range("A266").Value = 4
$ws.Range("B266").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C266").Value = "Los Lagos"
$ws.Range("D266").Value = 44943
$ws.Range("D266").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E266").Value = 10
$ws.Range("F266").Value = 100112044
$ws.Range("G266").Value = "Perejil"
$ws.Range("H266").Value = "Sin especificar"
$ws.Range("I266").Value = "Primera"
$ws.Range("J266").Value = 160
$ws.Range("K266").Value = 6000
$ws.Range("L266").Value = 6000
$ws.Range("M266").Value = 6000
$ws.Range("N266").Value = '$/docena de atados (2 kilos)'
$ws.Range("O266").Value = "Región de La Araucanía"
$ws.Range("P266").Value = 3000
$ws.Range("Q266").Value = 2
$ws.Range("R266").Value = "Hortaliza"
